$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds a date serial value for every data row (rows 2-278).
# All of them are being bumped by one day (46060 -> 46061).
$range = $ws.Range("C2:C278")
$range.Value = 46061
